# Update "想去人数" (people-interested count) figures across sheets to
# reflect the latest scraped totals (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3375
$ws.Range("F5").Value = 221
$ws.Range("F6").Value = 4870
$ws.Range("F7").Value = 475
$ws.Range("F9").Value = 180
$ws.Range("F10").Value = 637
$ws.Range("F11").Value = 287
$ws.Range("F12").Value = 46
$ws.Range("F14").Value = 667
$ws.Range("F19").Value = 349
$ws.Range("F20").Value = 4777
$ws.Range("F21").Value = 28
$ws.Range("F24").Value = 5914
$ws.Range("F25").Value = 17
$ws.Range("F26").Value = 1202
$ws.Range("F27").Value = 248
$ws.Range("F28").Value = 680
$ws.Range("F29").Value = 4427
$ws.Range("F30").Value = 8
$ws.Range("F33").Value = 882
$ws.Range("F35").Value = 11
$ws.Range("F36").Value = 800
$ws.Range("F37").Value = 871
$ws.Range("F38").Value = 3

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 15

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 222
$ws.Range("F4").Value = 42

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 222
$ws.Range("F5").Value = 42
$ws.Range("F8").Value = 3375
$ws.Range("F9").Value = 221
$ws.Range("F10").Value = 4870
$ws.Range("F11").Value = 475
$ws.Range("F13").Value = 180
$ws.Range("F14").Value = 637
$ws.Range("F15").Value = 287
$ws.Range("F16").Value = 46
$ws.Range("F18").Value = 667
$ws.Range("F24").Value = 349
$ws.Range("F25").Value = 4777
$ws.Range("F26").Value = 28
$ws.Range("F29").Value = 5914
$ws.Range("F30").Value = 17
$ws.Range("F31").Value = 1202
$ws.Range("F32").Value = 248
$ws.Range("F33").Value = 680
$ws.Range("F34").Value = 4427
$ws.Range("F35").Value = 8
$ws.Range("F36").Value = 15
$ws.Range("F39").Value = 882
$ws.Range("F41").Value = 11
$ws.Range("F42").Value = 800
$ws.Range("F43").Value = 871
$ws.Range("F45").Value = 3
